# Update to US commit ecc67274 on 6/5/24
#
# 1. Fill in the SoCiIEPTtB sheet with the 25-industry breakdown table
#    (replacing the single "Share of expenses" placeholder row), each row
#    referencing the full-passthrough assumption in About!A42.
# 2. Clean up the "About" sheet cosmetics that came along with the resave
#    (drop the now-unused italic style on the source-citation lines).
# 3. Move the active/selected sheet & cell to match the saved state.

$wb = $excel.ActiveWorkbook
$about = $wb.Worksheets.Item("About")
$data = $wb.Worksheets.Item("SoCiIEPTtB")

$industries = @(
    "agriculture and forestry 01T03",
    "coal mining 05",
    "oil and gas extraction 06",
    "other mining and quarrying 07T08",
    "food beverage and tobacco 10T12",
    "textiles apparel and leather 13T15",
    "wood products 16",
    "pulp paper and printing 17T18",
    "refined petroleum and coke 19",
    "chemicals 20",
    "rubber and plastic products 22",
    "glass and glass products 231",
    "cement and other nonmetallic minerals 239",
    "iron and steel 241",
    "other metals 242",
    "metal products except machinery and vehicles 25",
    "computers and electronics 26",
    "appliances and electrical equipment 27",
    "other machinery 28",
    "road vehicles 29",
    "nonroad vehicles 30",
    "other manufacturing 31T33",
    "energy pipelines and gas processing 352T353",
    "water and waste 36T39",
    "construction 41T43"
)

# Row 2 already exists ("Share of expenses" -> first industry); rows 3-26
# are brand new. Column B always repeats the full-passthrough flag from
# the About sheet, now pinned to the row with an absolute reference.
$row = 2
foreach ($industry in $industries) {
    $data.Cells.Item($row, 1).Value = $industry
    $data.Cells.Item($row, 2).Formula = "=About!A`$42"
    $row = $row + 1
}

# The two citation lines at the bottom of "About" lose their italic
# "Normal"-adjacent style in the resave.
$about.Cells.Item(40, 1).Style = "Normal"
$about.Cells.Item(41, 1).Style = "Normal"

# Selection / active sheet bookkeeping that Excel records on save.
$about.Range("B42").Select() | Out-Null
$data.Range("I16").Select() | Out-Null
$data.Activate() | Out-Null
